$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I8").Value = 'sd'
$ws.Range("J8").Value = 'Statement-non-opinion'
$ws.Range("I23").Value = 'sd'
$ws.Range("J23").Value = 'Statement-non-opinion'
$ws.Range("I27").Value = 'sd'
$ws.Range("J27").Value = 'Statement-non-opinion'
$ws.Range("I37").Value = 'sv'
$ws.Range("J37").Value = 'Statement-opinion'
$ws.Range("I44").Value = 'sd'
$ws.Range("J44").Value = 'Statement-non-opinion'
$ws.Range("I45").Value = 'qy'
$ws.Range("J45").Value = 'Yes-No-Question'
$ws.Range("I46").Value = 'ba'
$ws.Range("J46").Value = 'Appreciation'
$ws.Range("I48").Value = 'ba'
$ws.Range("J48").Value = 'Appreciation'
$ws.Range("I50").Value = 'ba'
$ws.Range("J50").Value = 'Appreciation'
$ws.Range("I51").Value = 'sd'
$ws.Range("J51").Value = 'Statement-non-opinion'
$ws.Range("I54").Value = 'sv'
$ws.Range("J54").Value = 'Statement-opinion'
$ws.Range("I57").Value = 'sv'
$ws.Range("J57").Value = 'Statement-opinion'
$ws.Range("I58").Value = 'sv'
$ws.Range("J58").Value = 'Statement-opinion'
$ws.Range("I59").Value = 'sv'
$ws.Range("J59").Value = 'Statement-opinion'
$ws.Range("I62").Value = 'sd'
$ws.Range("J62").Value = 'Statement-non-opinion'
$ws.Range("I71").Value = 'sd'
$ws.Range("J71").Value = 'Statement-non-opinion'
$ws.Range("I103").Value = 'sd'
$ws.Range("J103").Value = 'Statement-non-opinion'
$ws.Range("I110").Value = 'sd'
$ws.Range("J110").Value = 'Statement-non-opinion'
$ws.Range("I112").Value = 'sv'
$ws.Range("J112").Value = 'Statement-opinion'
$ws.Range("I114").Value = 'sv'
$ws.Range("J114").Value = 'Statement-opinion'
$ws.Range("I117").Value = 'aa'
$ws.Range("J117").Value = 'Agree/Accept'
$ws.Range("I119").Value = '%'
$ws.Range("J119").Value = 'Uninterpretable'
$ws.Range("I120").Value = '%'
$ws.Range("J120").Value = 'Uninterpretable'
$ws.Range("I122").Value = 'sd'
$ws.Range("J122").Value = 'Statement-non-opinion'
$ws.Range("I123").Value = 'sd'
$ws.Range("J123").Value = 'Statement-non-opinion'
$ws.Range("I126").Value = 'aa'
$ws.Range("J126").Value = 'Agree/Accept'
$ws.Range("I143").Value = 'ba'
$ws.Range("J143").Value = 'Appreciation'
$ws.Range("I158").Value = 'sd'
$ws.Range("J158").Value = 'Statement-non-opinion'
$ws.Range("I159").Value = 'sd'
$ws.Range("J159").Value = 'Statement-non-opinion'
$ws.Range("I162").Value = 'sd'
$ws.Range("J162").Value = 'Statement-non-opinion'
$ws.Range("I163").Value = 'sd'
$ws.Range("J163").Value = 'Statement-non-opinion'
$ws.Range("I166").Value = '%'
$ws.Range("J166").Value = 'Uninterpretable'
$ws.Range("I172").Value = 'qy'
$ws.Range("J172").Value = 'Yes-No-Question'
$ws.Range("I179").Value = 'sd'
$ws.Range("J179").Value = 'Statement-non-opinion'
$ws.Range("I181").Value = 'sd'
$ws.Range("J181").Value = 'Statement-non-opinion'
$ws.Range("I185").Value = 'sd'
$ws.Range("J185").Value = 'Statement-non-opinion'
$ws.Range("I193").Value = 'sd'
$ws.Range("J193").Value = 'Statement-non-opinion'
$ws.Range("I199").Value = 'b'
$ws.Range("J199").Value = 'Acknowledge (Backchannel)'
$ws.Range("I203").Value = 'sv'
$ws.Range("J203").Value = 'Statement-opinion'
$ws.Range("I215").Value = 'sv'
$ws.Range("J215").Value = 'Statement-opinion'
$ws.Range("I223").Value = 'sd'
$ws.Range("J223").Value = 'Statement-non-opinion'
$ws.Range("I232").Value = 'aa'
$ws.Range("J232").Value = 'Agree/Accept'
$ws.Range("I235").Value = 'sd'
$ws.Range("J235").Value = 'Statement-non-opinion'
$ws.Range("I250").Value = 'sv'
$ws.Range("J250").Value = 'Statement-opinion'
$ws.Range("I253").Value = 'sv'
$ws.Range("J253").Value = 'Statement-opinion'
$ws.Range("I268").Value = 'sd'
$ws.Range("J268").Value = 'Statement-non-opinion'
$ws.Range("I270").Value = 'sd'
$ws.Range("J270").Value = 'Statement-non-opinion'
$ws.Range("I284").Value = 'sv'
$ws.Range("J284").Value = 'Statement-opinion'
$ws.Range("I286").Value = 'sv'
$ws.Range("J286").Value = 'Statement-opinion'
$ws.Range("I296").Value = '%'
$ws.Range("J296").Value = 'Uninterpretable'
$ws.Range("I315").Value = 'sd'
$ws.Range("J315").Value = 'Statement-non-opinion'
$ws.Range("I320").Value = 'sv'
$ws.Range("J320").Value = 'Statement-opinion'
$ws.Range("I324").Value = 'aa'
$ws.Range("J324").Value = 'Agree/Accept'
$ws.Range("I325").Value = 'sv'
$ws.Range("J325").Value = 'Statement-opinion'
$ws.Range("I331").Value = 'sv'
$ws.Range("J331").Value = 'Statement-opinion'
